$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 52, shifting existing rows 52-122 down to 53-123
$ws.Rows.Item(52).Insert()

$ws.Cells.Item(52, 1).Value = 2
$ws.Cells.Item(52, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(52, 3).Value = "Coquimbo"
$ws.Cells.Item(52, 4).Value = 44665
$ws.Cells.Item(52, 5).Value = 4
$ws.Cells.Item(52, 6).Value = 100112024
$ws.Cells.Item(52, 7).Value = "Choclo"
$ws.Cells.Item(52, 8).Value = "Choclero"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 40000
$ws.Cells.Item(52, 11).Value = 230
$ws.Cells.Item(52, 12).Value = 250
$ws.Cells.Item(52, 13).Value = 240
$ws.Cells.Item(52, 14).Value = "$/unidad"
$ws.Cells.Item(52, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(52, 16).Value = 240
$ws.Cells.Item(52, 17).Value = 1
$ws.Cells.Item(52, 18).Value = "Hortaliza"
